$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Periodo Mora" (column E) and "Valor Mora" (column F) for rows 16 and 18
# (rows 16 and 18 values are effectively swapped)
$ws.Range("E16").Value = "2410"
$ws.Range("F16").Value = 52000

$ws.Range("E18").Value = "2502"
$ws.Range("F18").Value = 56940

# Update "Salario Basico" (column G) for rows 16, 17 and 18 to the new value
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
